# Apply the change described by the diff:
#  - insert a new column E ("neighbourhood_group") before the old
#    "Bydel_area" column, shifting Bydel_area/longitude/latitude from
#    E/F/G to F/G/H
#  - populate the new column's header and the 16 data values
#  - fix D6's text from "St.Hanshaugen" to "St. Hanshaugen"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting existing E:G to F:H.
$ws.Range("E1").EntireColumn.Insert()

# New header.
$ws.Range("E1").Value = "neighbourhood_group"

# New column values (neighbourhood_group), one per data row.
$ws.Range("E2").Value = 16
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 9
$ws.Range("E6").Value = 13
$ws.Range("E7").Value = 7
$ws.Range("E8").Value = 10
$ws.Range("E9").Value = 5
$ws.Range("E10").Value = 11
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 4
$ws.Range("E13").Value = 14
$ws.Range("E14").Value = 1
$ws.Range("E15").Value = 15
$ws.Range("E16").Value = 8
$ws.Range("E17").Value = 12

# Fix a typo in the Bydel_Navn column.
$ws.Range("D6").Value = "St. Hanshaugen"
